# Deploy the implementation guide: refresh the "observation-code" CodeSystem
# export (Title/Date/Count metadata + the full Concepts table).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Concepts")

# --- Metadata sheet -------------------------------------------------------
$ws1.Range("B5").Value  = "Observation codes"          # Title
$ws1.Range("B8").Value  = "2022-05-18T17:38:26+00:00"  # Date

# "Count" (B21) must stay text like the rest of the column even though "25"
# looks numeric. A plain .Value assignment on a General-formatted cell gets
# auto-coerced to a number, so build the literal text via a scratch-cell
# formula and paste its *value* back in (mirrors typing '25 in the UI
# without leaving a quote-prefix / number-format behind).
$scratch = $ws1.Range("Z1")
$scratch.Formula = "=""25"""
$scratch.Copy()
$ws1.Range("B21").PasteSpecial(-4163)
$scratch.ClearContents()

# --- Concepts sheet ---------------------------------------------------------
# New full set of concepts (Level is always "1"; Definition column is blank).
$data = @(
  @('OBSG', 'General Observation'),
  @('INDIC', 'Indication'),
  @('INVES', 'Investigation'),
  @('BMUS', 'Muscle biopsy'),
  @('BMET', 'Metabolic work-up'),
  @('CKIN', 'Serum creatine kinase'),
  @('CAAP', 'Plasma amino acid chromatography'),
  @('ACYL', 'Acylcarnitines'),
  @('IRMC', 'Cerebral MRI'),
  @('IRMM', 'Muscle MRI'),
  @('SRMN', 'NMR Spectroscopy'),
  @('ECAR', 'Cardiac Ultrasound'),
  @('EABD', 'Abdominal Ultrasound'),
  @('EEG', 'Electroencephalogram (EEG)'),
  @('EMG', 'Electromyography (EMG)'),
  @('OPHT', 'Ophthalmological Examination'),
  @('CONS', 'Consanguinity'),
  @('DSTA', 'Patient Disease Status'),
  @('CTGR', 'CTG Repeat Testing (Steinert Disease)'),
  @('GCNR', 'GCN Repeat Testing (oculopharyngeal muscular dystrophy)'),
  @('TGDD', 'Deletions and Duplications Testing (Duchenne and Becker dystrophies)'),
  @('TSXF', 'Fragile X Syndrome Testing'),
  @('MFTH', 'Missing Mother'),
  @('MMTH', 'Missing Father'),
  @('ETHN', 'Ethnicity')
)

$lastExistingRow = 7     # rows 2..7 already exist in the sheet
$firstRow = 2

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $firstRow + $i
  if ($r -gt $lastExistingRow) {
    # Brand-new row: clone the formatting (incl. the text-typed "1" in
    # column A) from the row above before writing the new B/C values.
    $ws2.Range("A$($r-1):D$($r-1)").Copy($ws2.Range("A$($r):D$($r)"))
  }
  $ws2.Cells.Item($r, 2).Value = $data[$i][0]
  $ws2.Cells.Item($r, 3).Value = $data[$i][1]
  $ws2.Cells.Item($r, 4).ClearContents()
}

Write-Host "done"
